# Fit or Fail_presentation.pptx — "Add files via upload" edit
#
# 1. The slide-master and every slide-layout carry a cached
#    datetimeFigureOut field ("3/4/2020" -> "3/10/2020").
# 2. Slide 7's status table: the "Use Cases/Requirements" row's
#    "Who Completed" cell changes from "Christian" to "Josh/Christian".

$p = $ppt.ActivePresentation

# --- 1. Update the cached date text on the master + every layout -----------

$newDate = "3/10/2020"

$master = $p.SlideMaster
for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $shape = $master.Shapes.Item($si)
    if ($shape.Name -like "Date Placeholder*") {
        $shape.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shape = $layout.Shapes.Item($si)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Slide 7 status table: "Christian" -> "Josh/Christian" --------------

$slide7 = $p.Slides.Item(7)
for ($si = 1; $si -le $slide7.Shapes.Count; $si++) {
    $shape = $slide7.Shapes.Item($si)
    if ($shape.HasTable) {
        $table = $shape.Table
        for ($r = 1; $r -le $table.Rows.Count; $r++) {
            for ($c = 1; $c -le $table.Columns.Count; $c++) {
                $cell = $table.Cell($r, $c)
                $tr = $cell.Shape.TextFrame.TextRange
                if ($tr.Text -eq "Christian") {
                    $tr.Text = "Josh/Christian"
                }
            }
        }
    }
}
